$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (unmodified, default-formatted cell) used to restore
# the default number format on text cells after forcing text entry,
# so numeric-looking strings (e.g. "553.27") stay text without leaving
# a stray explicit style on the cell.
$refStyle = $ws.Range("B2").Style

function Set-TextValue([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $refStyle
}

$ws.Range("D2").Value = '62.821.44'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '2.680.31'
$ws.Range("E3").Value = '  -2.01%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '553.27'
$ws.Range("E5").Value = '  -2.17%  '
Set-TextValue "D6" '158.12'
$ws.Range("E6").Value = '  -0.72%  '
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  +0.06%  '
Set-TextValue "D8" '0.590'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("E9").Value = '  -2.64%  '
$ws.Range("E10").Value = '  -2.19%  '
Set-TextValue "D11" '0.369'
$ws.Range("E11").Value = '  -2.99%  '
Set-TextValue "D12" '5.40'
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("D13").Value = '3.155.81'
$ws.Range("E13").Value = '  -1.94%  '
Set-TextValue "D14" '26.50'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").Value = '62.745.65'
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '2.682.14'
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("E18").Value = '  -3.82%  '
$ws.Range("E19").Value = '  -2.66%  '
Set-TextValue "D20" '344.90'
$ws.Range("E20").Value = '  -2.33%  '
$ws.Range("E21").Value = '  -4.25%  '
$ws.Range("E22").Value = '  +0.00%  '
Set-TextValue "D23" '0.507'
$ws.Range("E23").Value = '  -2.88%  '
Set-TextValue "D24" '63.11'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -0.03%  '
Set-TextValue "D27" '8.19'
$ws.Range("E27").Value = '  -1.80%  '
$ws.Range("E28").Value = '  +9.29%  '
$ws.Range("D29").Value = '0.0₃0855'
$ws.Range("E29").Value = '  -5.22%  '
Set-TextValue "D30" '7.24'
$ws.Range("E30").Value = '  +0.75%  '
Set-TextValue "D31" '1.94'
$ws.Range("E31").Value = '  -1.09%  '
Set-TextValue "D32" '163.99'
$ws.Range("E32").Value = '  +0.36%  '
Set-TextValue "D33" '4.93'
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D34" '1.48'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D35" '0.999'
$ws.Range("E35").Value = '  +0.00%  '
Set-TextValue "D36" '19.48'
$ws.Range("E36").Value = '  -2.74%  '
$ws.Range("E37").Value = '  -0.06%  '
Set-TextValue "D38" '348.75'
$ws.Range("E38").Value = '  +1.32%  '
Set-TextValue "D39" '0.947'
$ws.Range("E39").Value = '  -3.26%  '
Set-TextValue "D40" '6.23'
$ws.Range("E40").Value = '  -0.74%  '
Set-TextValue "D41" '3.98'
$ws.Range("E41").Value = '  -2.01%  '
Set-TextValue "D42" '38.40'
$ws.Range("E42").Value = '  -0.07%  '
Set-TextValue "D43" '20.92'
$ws.Range("E43").Value = '  -3.71%  '
Set-TextValue "D44" '20.16'
$ws.Range("E44").Value = '  -3.55%  '
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("E47").Value = '  -0.01%  '
Set-TextValue "D48" '11.02'
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D49" '0.0970'
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D50" '0.0242'
$ws.Range("E50").Value = '  -3.04%  '
Set-TextValue "D51" '129.02'
$ws.Range("E51").Value = '  -3.86%  '
